$d = $word.ActiveDocument

# Update the date line
$d.Content.Find.Execute("2024-03-20 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-03-21 Thursday", 2) | Out-Null

# Update table cell values (table is 20 rows x 5 columns)
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "97-13=84"
$t.Cell(1,2).Range.Text = "82-16=66"
$t.Cell(1,3).Range.Text = "18-14=4"
$t.Cell(1,4).Range.Text = "77+11=88"
$t.Cell(1,5).Range.Text = "68+8=76"

$t.Cell(2,1).Range.Text = "6+34=40"
$t.Cell(2,2).Range.Text = "51-15=36"
$t.Cell(2,3).Range.Text = "98-93=5"
$t.Cell(2,4).Range.Text = "19+45=64"
$t.Cell(2,5).Range.Text = "45-44=1"

$t.Cell(3,1).Range.Text = "79-0=79"
$t.Cell(3,2).Range.Text = "39+10=49"
$t.Cell(3,3).Range.Text = "38+43=81"
$t.Cell(3,4).Range.Text = "4+68=72"
$t.Cell(3,5).Range.Text = "0+91=91"

$t.Cell(4,1).Range.Text = "25+73=98"
$t.Cell(4,2).Range.Text = "43+43=86"
$t.Cell(4,3).Range.Text = "63-56=7"
$t.Cell(4,4).Range.Text = "27-3=24"
$t.Cell(4,5).Range.Text = "35+21=56"

$t.Cell(5,1).Range.Text = "71+8=79"
$t.Cell(5,2).Range.Text = "8+53=61"
$t.Cell(5,3).Range.Text = "95-40=55"
$t.Cell(5,4).Range.Text = "33+15=48"
$t.Cell(5,5).Range.Text = "66+7=73"

$t.Cell(6,1).Range.Text = "80-42=38"
$t.Cell(6,2).Range.Text = "51+15=66"
$t.Cell(6,3).Range.Text = "92-33=59"
$t.Cell(6,4).Range.Text = "36+34=70"
$t.Cell(6,5).Range.Text = "8+14=22"

$t.Cell(7,1).Range.Text = "17+23=40"
$t.Cell(7,2).Range.Text = "29-18=11"
$t.Cell(7,3).Range.Text = "36+56=92"
$t.Cell(7,4).Range.Text = "23+21=44"
$t.Cell(7,5).Range.Text = "15+11=26"

$t.Cell(8,1).Range.Text = "41-3=38"
$t.Cell(8,2).Range.Text = "46-23=23"
$t.Cell(8,3).Range.Text = "65+33=98"
$t.Cell(8,4).Range.Text = "2+54=56"
$t.Cell(8,5).Range.Text = "20+56=76"

$t.Cell(9,1).Range.Text = "73-70=3"
$t.Cell(9,2).Range.Text = "93-8=85"
$t.Cell(9,3).Range.Text = "50+41=91"
$t.Cell(9,4).Range.Text = "65+18=83"
$t.Cell(9,5).Range.Text = "52+33=85"

$t.Cell(10,1).Range.Text = "38+17=55"
$t.Cell(10,2).Range.Text = "53-20=33"
$t.Cell(10,3).Range.Text = "35+64=99"
$t.Cell(10,4).Range.Text = "44+44=88"
$t.Cell(10,5).Range.Text = "86-1=85"

$t.Cell(11,1).Range.Text = "46+7=53"
$t.Cell(11,2).Range.Text = "83-55=28"
$t.Cell(11,3).Range.Text = "25-23=2"
$t.Cell(11,4).Range.Text = "72-10=62"
$t.Cell(11,5).Range.Text = "62-1=61"

$t.Cell(12,1).Range.Text = "61-23=38"
$t.Cell(12,2).Range.Text = "77+9=86"
$t.Cell(12,3).Range.Text = "99-43=56"
$t.Cell(12,4).Range.Text = "11+22=33"
$t.Cell(12,5).Range.Text = "3+54=57"

$t.Cell(13,1).Range.Text = "97-75=22"
$t.Cell(13,2).Range.Text = "42-25=17"
$t.Cell(13,3).Range.Text = "53-27=26"
$t.Cell(13,4).Range.Text = "10+73=83"
$t.Cell(13,5).Range.Text = "61-39=22"

$t.Cell(14,1).Range.Text = "8+36=44"
$t.Cell(14,2).Range.Text = "1+12=13"
$t.Cell(14,3).Range.Text = "12-4=8"
$t.Cell(14,4).Range.Text = "74-19=55"
$t.Cell(14,5).Range.Text = "7+66=73"

$t.Cell(15,1).Range.Text = "5+19=24"
$t.Cell(15,2).Range.Text = "99-42=57"
$t.Cell(15,3).Range.Text = "98-25=73"
$t.Cell(15,4).Range.Text = "85-11=74"
$t.Cell(15,5).Range.Text = "48+26=74"

$t.Cell(16,1).Range.Text = "81+5=86"
$t.Cell(16,2).Range.Text = "57-34=23"
$t.Cell(16,3).Range.Text = "46-8=38"
$t.Cell(16,4).Range.Text = "60+24=84"
$t.Cell(16,5).Range.Text = "64-5=59"

$t.Cell(17,1).Range.Text = "71-68=3"
$t.Cell(17,2).Range.Text = "51-1=50"
$t.Cell(17,3).Range.Text = "74-45=29"
$t.Cell(17,4).Range.Text = "81+5=86"
$t.Cell(17,5).Range.Text = "58+22=80"

$t.Cell(18,1).Range.Text = "74-67=7"
$t.Cell(18,2).Range.Text = "43+54=97"
$t.Cell(18,3).Range.Text = "53+41=94"
$t.Cell(18,4).Range.Text = "4+49=53"
$t.Cell(18,5).Range.Text = "51+9=60"

$t.Cell(19,1).Range.Text = "8-4=4"
$t.Cell(19,2).Range.Text = "58+23=81"
$t.Cell(19,3).Range.Text = "97-35=62"
$t.Cell(19,4).Range.Text = "60-37=23"
$t.Cell(19,5).Range.Text = "62+36=98"

$t.Cell(20,1).Range.Text = "79-20=59"
$t.Cell(20,2).Range.Text = "50-22=28"
$t.Cell(20,3).Range.Text = "48-10=38"
$t.Cell(20,4).Range.Text = "90+7=97"
$t.Cell(20,5).Range.Text = "13+40=53"

Write-Host "Done updating cells"